$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.113.75"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.667.08"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'209.80"
$ws.Range("D6").Value = "'0.5209"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D8").Value = "'0.2598"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "'0.06324"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'21.06"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.07529"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.672.06"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "'4.410"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "'0.5411"
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "'0.000007995"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "26.190.86"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "'186.99"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("D22").Value = "'6.221"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'149.73"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "'0.1238"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "'7.440"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").Value = "'15.72"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "'0.06284"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").Value = "'1.275"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "'3.494"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").Value = "'1.635"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "'0.9992"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.395"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.762"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.5974"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "1.110.20"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").Value = "'0.01609"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'6.056"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "'0.8615"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D43").Value = "'100.64"
$ws.Range("D44").Value = "1.815.89"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Value = "'55.33"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "'8.052"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "'0.05253"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'0.4234"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").Value = "'5.878"
$ws.Range("E51").Value = "  -1.08%  "
